$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tarefas")

# New row of task data (row 59), following the same pattern as the rows above.
$ws.Cells.Item(59, 1).Value = Get-Date -Year 2015 -Month 7 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(59, 2).Value = "Em Andamento"
$ws.Cells.Item(59, 3).Value = "Codificação"
$ws.Cells.Item(59, 4).Value = "OverFlow quando usa a função Clone"

# Copy formatting from the row above so the new row matches the sheet's look.
$ws.Range("A58:D58").Copy() | Out-Null
$ws.Range("A59:D59").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection to match where the user ended up after typing.
$ws.Range("D60").Select() | Out-Null
